# Generate Report for Handoff
# Adds a new handed-off file (07cec6c5-5fda-43d3-be41-0e222aa8f8ed) ahead of
# the existing 4a79108b-... entry on every sheet, and pushes the static
# ".localization-config" row down by one.

$wb = $excel.ActiveWorkbook

function Set-HyperlinkCell {
    param($ws, $cellRef, $target, $display)
    $ws.Hyperlinks.Add($ws.Range($cellRef), $target, "", "", $display)
}

# ---------------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de
# ---------------------------------------------------------------------------
$ov = $wb.Worksheets.Item("Overview")

# Drop the existing hyperlinks so they can be rebuilt in final left-to-right,
# top-to-bottom order (keeps relationship ids / hyperlink ordering sane).
$ov.Hyperlinks.Delete()

$ov.Range("A2").Value = "07cec6c5-5fda-43d3-be41-0e222aa8f8ed.md"
$ov.Range("B2").Value = "Ready for handoff"
$ov.Range("C2").Value = "Ready for handoff"

$ov.Range("A3").Value = "4a79108b-b4a0-479e-a732-f05fc04daeae.md"
$ov.Range("B3").Value = "Ready for handoff"
$ov.Range("C3").Value = "Ready for handoff"

$ov.Range("A4").Value = ".localization-config"
$ov.Range("B4").Value = "Not to be localized"
$ov.Range("C4").Value = "Not to be localized"

$ov.Range("A2:A4").Font.Underline = $true
$ov.Range("A2:A4").Font.Color = 13011546

Set-HyperlinkCell $ov "A2" "https://github.com/OpenLocalizationTest/oltest/blob/e2e/07cec6c5-5fda-43d3-be41-0e222aa8f8ed.md" "07cec6c5-5fda-43d3-be41-0e222aa8f8ed.md"
Set-HyperlinkCell $ov "A3" "https://github.com/OpenLocalizationTest/oltest/blob/e2e/4a79108b-b4a0-479e-a732-f05fc04daeae.md" "4a79108b-b4a0-479e-a732-f05fc04daeae.md"
Set-HyperlinkCell $ov "A4" "https://github.com/OpenLocalizationTest/oltest/blob/.localization-config" ".localization-config"

$ov.Range("A1:C4").Columns.AutoFit() | Out-Null

# ---------------------------------------------------------------------------
# Per-language detail sheets: Source File Name | Status | Latest Handoff File
# | Latest Handoff Datetime | Latest Target File | Latest Handback File |
# Latest Handback DateTime | Handoff Reason | Dependency From
# ---------------------------------------------------------------------------
function Update-LangSheet {
    param($ws, $newXlf, $newDate, $oldXlf, $oldDate)

    $ws.Hyperlinks.Delete()

    # Row 2: newly generated handoff file
    $ws.Range("A2").Value = "07cec6c5-5fda-43d3-be41-0e222aa8f8ed.md"
    $ws.Range("B2").Value = "Ready for handoff"
    $ws.Range("C2").Value = $newXlf
    $ws.Range("D2").Value = $newDate
    $ws.Range("E2").Value = ""
    $ws.Range("F2").Value = ""
    $ws.Range("G2").Value = "0001-01-01 00:00:00"
    $ws.Range("H2").Value = "Include"
    $ws.Range("I2").Value = ""

    # Row 3: previous handoff file (was row 2 before this run)
    $ws.Range("A3").Value = "4a79108b-b4a0-479e-a732-f05fc04daeae.md"
    $ws.Range("B3").Value = "Ready for handoff"
    $ws.Range("C3").Value = $oldXlf
    $ws.Range("D3").Value = $oldDate
    $ws.Range("E3").Value = ""
    $ws.Range("F3").Value = ""
    $ws.Range("G3").Value = "0001-01-01 00:00:00"
    $ws.Range("H3").Value = "Include"
    $ws.Range("I3").Value = ""

    # Row 4: static ".localization-config" row (was row 3 before this run)
    $ws.Range("A4").Value = ".localization-config"
    $ws.Range("B4").Value = "Not to be localized"
    $ws.Range("C4").Value = ""
    $ws.Range("D4").Value = "0001-01-01 00:00:00"
    $ws.Range("E4").Value = ""
    $ws.Range("F4").Value = ""
    $ws.Range("G4").Value = "0001-01-01 00:00:00"
    $ws.Range("H4").Value = "Ignored"
    $ws.Range("I4").Value = ""

    $ws.Range("A2:A4").Font.Underline = $true
    $ws.Range("A2:A4").Font.Color = 13011546
    $ws.Range("C2:C3").Font.Underline = $true
    $ws.Range("C2:C3").Font.Color = 13011546
    $ws.Range("D2:D4").NumberFormat = "yyyy-mm-dd HH:mm:ss"

    Set-HyperlinkCell $ws "A2" "https://github.com/OpenLocalizationTest/oltest/blob/e2e/07cec6c5-5fda-43d3-be41-0e222aa8f8ed.md" "07cec6c5-5fda-43d3-be41-0e222aa8f8ed.md"
    Set-HyperlinkCell $ws "C2" ("https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ht/" + $newXlf) $newXlf
    Set-HyperlinkCell $ws "A3" "https://github.com/OpenLocalizationTest/oltest/blob/e2e/4a79108b-b4a0-479e-a732-f05fc04daeae.md" "4a79108b-b4a0-479e-a732-f05fc04daeae.md"
    Set-HyperlinkCell $ws "C3" ("https://github.com/OpenLocalizationTestOrg/olhandoff/blob/ht/" + $oldXlf) $oldXlf
    Set-HyperlinkCell $ws "A4" "https://github.com/OpenLocalizationTest/oltest/blob/.localization-config" ".localization-config"

    $ws.Range("A1:I4").Columns.AutoFit() | Out-Null
}

$zhcn = $wb.Worksheets.Item("zh-cn")
Update-LangSheet $zhcn `
    "07cec6c5-5fda-43d3-be41-0e222aa8f8ed.4989a88dc5ae437817d156855b16b3dd4195a646.zh-cn.xlf" "2016-03-10 22:43:58" `
    "4a79108b-b4a0-479e-a732-f05fc04daeae.a66a1967ac86ae38aec6e8788dd7635b7e5215f7.zh-cn.xlf" "2016-03-10 22:43:31"

$dede = $wb.Worksheets.Item("de-de")
Update-LangSheet $dede `
    "07cec6c5-5fda-43d3-be41-0e222aa8f8ed.4989a88dc5ae437817d156855b16b3dd4195a646.de-de.xlf" "2016-03-10 22:44:05" `
    "4a79108b-b4a0-479e-a732-f05fc04daeae.a66a1967ac86ae38aec6e8788dd7635b7e5215f7.de-de.xlf" "2016-03-10 22:43:37"

Write-Output "Report regenerated for handoff"
